# Apply the changes described by the diff:
#  1. Update the weight (G20) and count (H20) for the Parapenaeus longirostris
#     entry on row 20 (same species as catch group "1-RAP").
#  2. Remove the duplicate "Parapenaeus longirostris / PAPELON" row (row 38,
#     catch group "2-RAP") which shifts all the following rows (39-44) up by
#     one, so the former rows 39-44 become rows 38-43 and the sheet ends up
#     with one fewer row (old dimension A1:K44 -> new dimension A1:K43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update row 20 values.
$ws.Range("G20").Value = 0.33
$ws.Range("H20").Value = 42

# 2. Delete row 38 entirely; Excel shifts rows 39+ up automatically and the
#    worksheet dimension is recalculated to A1:K43.
$ws.Rows.Item(38).Delete()
